$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Amount in Pool" values for the resource rows.
# Standard_B2ms (row 3): 3 -> 1
$ws.Range("E3").Value = 1

# Standard_B4ms (row 5): 0 -> 1
$ws.Range("E5").Value = 1

$wb.Save()
